# Scraper refresh: "30/12 21:50 LP1912+6203+6173"
# New rows appended to each of the 3 sheets with the scrape that ran at
# 18:50:xx, plus the "Última actualización" / "Total filas" header cells
# on each sheet bumped to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (cols: A blank-marker, B Hora_Scrap, C Hora_Llegada,
#                    D Linea, E Minutos, F Parada, G Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 30/12/2025 18:50:55"
$ws1.Range("A3").Value = "Total filas: 542"

$rows1 = @(
    @("18:50:44", "18:52", "15_ABASTO", 2),
    @("18:50:44", "18:56", "10_OLMOS", 6),
    @("18:50:44", "19:00", "16_SANTA ANA", 10),
    @("18:50:44", "19:04", "11_ETCHEVERRY", 14),
    @("18:50:44", "19:08", "23_HERNANDEZ", 18),
    @("18:50:44", "19:12", "10_OLMOS", 22),
    @("18:50:44", "19:20", "14_ABASTO", 30),
    @("18:50:44", "19:21", "16_SANTA ANA", 31),
    @("18:50:44", "19:21", "26_HERNANDEZ", 31),
    @("18:50:44", "19:28", "15_ABASTO", 38),
    @("18:50:44", "19:39", "215C_EL PATO", 49),
    @("18:50:44", "19:40", "14_ABASTO", 50),
    @("18:50:44", "19:50", "11X44_ETCHEVERRY", 60),
    @("18:50:44", "19:50", "16_P MOR-SANTA ANA", 60),
    @("18:50:44", "19:51", "81_EL PELIGRO", 61),
    @("18:50:44", "19:59", "17_ROMERO", 69),
    @("18:50:44", "20:00", "14_ABASTO", 70),
    @("18:50:44", "20:07", "10_OLMOS", 77),
    @("18:50:44", "20:09", "15_ABASTO", 79),
    @("18:50:44", "20:10", "16_P MOR-167 Y 521", 80),
    @("18:50:44", "20:12", "23_HERNANDEZ", 82),
    @("18:50:44", "20:20", "26_HERNANDEZ", 90),
    @("18:50:44", "20:22", "11_ETCHEVERRY", 92),
    @("18:50:44", "20:23", "215A_EL PATO", 93)
)

$r = 520
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = "LP1912"
    $ws1.Cells.Item($r, 7).Value = "30/12/2025"
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (cols: A blank-marker, B Fecha, C Hora_Scrap,
#                        D Hora_Llegada, E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 30/12/2025 18:50:55"
$ws2.Range("A3").Value = "Total filas: 35"

$rows2 = @(
    @("18:50:44", "19:39", "215C_EL PATO", 49),
    @("18:50:44", "20:23", "215A_EL PATO", 93)
)

$r = 35
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 2).Value = "30/12/2025"
    $ws2.Cells.Item($r, 3).Value = $row[0]
    $ws2.Cells.Item($r, 4).Value = $row[1]
    $ws2.Cells.Item($r, 5).Value = $row[2]
    $ws2.Cells.Item($r, 6).Value = $row[3]
    $ws2.Cells.Item($r, 7).Value = "LP1912"
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (cols: A blank-marker, B Fecha, C Hora_Scrap,
#                       D Hora_Llegada, E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 30/12/2025 18:50:55"
$ws3.Range("A3").Value = "Total filas: 70"

$rows3 = @(
    @("18:50:55", "19:03", "215B_LP-P MOR-1 Y 57", 13, "L6173"),
    @("18:50:50", "19:53", "215C_LA PLATA", 63, "L6203")
)

$r = 70
foreach ($row in $rows3) {
    $ws3.Cells.Item($r, 2).Value = "30/12/2025"
    $ws3.Cells.Item($r, 3).Value = $row[0]
    $ws3.Cells.Item($r, 4).Value = $row[1]
    $ws3.Cells.Item($r, 5).Value = $row[2]
    $ws3.Cells.Item($r, 6).Value = $row[3]
    $ws3.Cells.Item($r, 7).Value = $row[4]
    $r++
}
